$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: finish the time range and log hours / work description
$ws.Range("C6").Value = "12:30pm-3:30pm"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "Minor Changes to GUI, Code Review, Class Meeting, Team Meeting"

# Row 7: new day entry - 05/14/2015 Thu
$ws.Range("B7").Value = "05/14/2015 Thu"
$ws.Range("C7").Value = "5:45pm-9:45pm"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = "Search Box Formatting, Hash Function, Additional GUI Work, Code Cleanup"

# Row 8: new day entry - 05/15/2015 Fri (time range started, not yet finished)
$ws.Range("B8").Value = "05/15/2015 Fri"
$ws.Range("C8").Value = "12:00pm-"

# Widen the "Worked On" column to fit the newly added descriptions
$ws.Columns.Item(5).ColumnWidth = 68

# Move the active selection to the last edited cell
$ws.Range("C8").Select()
